$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("具有相當價值之財產")

# ---- Row 1: convert to a proper header row (like the other sheets) ----
$ws.Cells.Item(1,2).Value = "name"
$ws.Cells.Item(1,3).Value = "quantity"
$ws.Cells.Item(1,4).Value = "owner"
$ws.Cells.Item(1,5).Value = "total"

$headerCols = @(6,7,8,9,10,11,12)
$headerVals = @("property_category","category","date","legislator_name","legislator_id","source_file","index")
for ($i = 0; $i -lt $headerCols.Length; $i++) {
    $c = $headerCols[$i]
    $ws.Cells.Item(1,$c).Value = $headerVals[$i]
    # Match the existing header formatting used by B1:E1 (bold, bordered, centered)
    $ws.Cells.Item(1,$c).Font.Bold = $true
    $ws.Cells.Item(1,$c).HorizontalAlignment = -4108
    $ws.Cells.Item(1,$c).VerticalAlignment = -4160
    $ws.Cells.Item(1,$c).Borders.LineStyle = 1
}

# ---- Rows 2-6: append the metadata columns (F-L) that every other sheet already has ----
$indices = @(104,105,106,107,108)
for ($i = 0; $i -lt $indices.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row,6).Value  = "otherbonds"
    $ws.Cells.Item($row,7).Value  = "normal"
    # Force text so the date-like string isn't auto-converted to a date serial
    $ws.Cells.Item($row,8).NumberFormat = "@"
    $ws.Cells.Item($row,8).Value  = "2012-03-03"
    $ws.Cells.Item($row,9).Value  = "孫大千"
    $ws.Cells.Item($row,10).Value = 919
    $ws.Cells.Item($row,11).Value = "tmpc261"
    $ws.Cells.Item($row,12).Value = $indices[$i]
}
